$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("beach"): the is_hot/onoff column C7 becomes a real Boolean TRUE
# (previously stored as the plain number 1).
$ws.Range("C7").Value = $true

# Row 8 ("wonderland"): the is_hot/onoff column C8 becomes the literal text
# "FaLsE" (previously stored as the plain number 0). Format the cell as Text
# first, write a formula that evaluates to the text, then flatten the
# formula to a static value via copy / paste-values so it lands as a plain
# shared-string text cell (not a boolean, not a quote-prefixed literal).
$c8 = $ws.Range("C8")
$c8.NumberFormat = "@"
$c8.Formula = "=""FaLsE"""
$c8.Copy()
$c8.PasteSpecial(-4163)   # xlPasteValues

# Selection moves from A13 to C9.
$ws.Range("C9").Select() | Out-Null
